$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("re_potentials")

# Columns A-F (Sets, process, description, TAct, TCap, timeslicelevel)
$ws.Cells.Item(2,1).Value = 'Sets'
$ws.Cells.Item(2,2).Value = 'process'
$ws.Cells.Item(2,3).Value = 'description'
$ws.Cells.Item(2,4).Value = 'TAct'
$ws.Cells.Item(2,5).Value = 'TCap'
$ws.Cells.Item(2,6).Value = 'timeslicelevel'

$ws.Cells.Item(3,1).Value = 'ELE'
$ws.Cells.Item(3,2).Value = 'EN_SPV_13_c02_CHE'
$ws.Cells.Item(3,3).Value = 'Utility PV - CF Class-13 Cost Class-c02 - Switzerland'
$ws.Cells.Item(3,4).Value = 'TWh'
$ws.Cells.Item(3,5).Value = 'GW'
$ws.Cells.Item(3,6).Value = 'ANNUAL'

$ws.Cells.Item(4,1).Value = 'ELE'
$ws.Cells.Item(4,2).Value = 'EN_SPV_13_c03_CHE'
$ws.Cells.Item(4,3).Value = 'Utility PV - CF Class-13 Cost Class-c03 - Switzerland'
$ws.Cells.Item(4,4).Value = 'TWh'
$ws.Cells.Item(4,5).Value = 'GW'
$ws.Cells.Item(4,6).Value = 'ANNUAL'

$ws.Cells.Item(5,1).Value = 'ELE'
$ws.Cells.Item(5,2).Value = 'EN_SPV_14_c02_CHE'
$ws.Cells.Item(5,3).Value = 'Utility PV - CF Class-14 Cost Class-c02 - Switzerland'
$ws.Cells.Item(5,4).Value = 'TWh'
$ws.Cells.Item(5,5).Value = 'GW'
$ws.Cells.Item(5,6).Value = 'ANNUAL'

$ws.Cells.Item(6,1).Value = 'ELE'
$ws.Cells.Item(6,2).Value = 'EN_SPV_15_c02_CHE'
$ws.Cells.Item(6,3).Value = 'Utility PV - CF Class-15 Cost Class-c02 - Switzerland'
$ws.Cells.Item(6,4).Value = 'TWh'
$ws.Cells.Item(6,5).Value = 'GW'
$ws.Cells.Item(6,6).Value = 'ANNUAL'

$ws.Cells.Item(7,1).Value = 'ELE'
$ws.Cells.Item(7,2).Value = 'EN_WON_18_c04_CHE'
$ws.Cells.Item(7,3).Value = 'Wind Onshore - CF Class-18 Cost Class-c04 - Switzerland'
$ws.Cells.Item(7,4).Value = 'TWh'
$ws.Cells.Item(7,5).Value = 'GW'
$ws.Cells.Item(7,6).Value = 'ANNUAL'

$ws.Cells.Item(8,1).Value = 'ELE'
$ws.Cells.Item(8,2).Value = 'EN_WON_20_c03_CHE'
$ws.Cells.Item(8,3).Value = 'Wind Onshore - CF Class-20 Cost Class-c03 - Switzerland'
$ws.Cells.Item(8,4).Value = 'TWh'
$ws.Cells.Item(8,5).Value = 'GW'
$ws.Cells.Item(8,6).Value = 'ANNUAL'

$ws.Cells.Item(9,1).Value = 'ELE'
$ws.Cells.Item(9,2).Value = 'EN_WON_20_c04_CHE'
$ws.Cells.Item(9,3).Value = 'Wind Onshore - CF Class-20 Cost Class-c04 - Switzerland'
$ws.Cells.Item(9,4).Value = 'TWh'
$ws.Cells.Item(9,5).Value = 'GW'
$ws.Cells.Item(9,6).Value = 'ANNUAL'

$ws.Cells.Item(10,1).Value = 'ELE'
$ws.Cells.Item(10,2).Value = 'EN_WON_22_c02_CHE'
$ws.Cells.Item(10,3).Value = 'Wind Onshore - CF Class-22 Cost Class-c02 - Switzerland'
$ws.Cells.Item(10,4).Value = 'TWh'
$ws.Cells.Item(10,5).Value = 'GW'
$ws.Cells.Item(10,6).Value = 'ANNUAL'

$ws.Cells.Item(11,1).Value = 'ELE'
$ws.Cells.Item(11,2).Value = 'EN_WON_22_c03_CHE'
$ws.Cells.Item(11,3).Value = 'Wind Onshore - CF Class-22 Cost Class-c03 - Switzerland'
$ws.Cells.Item(11,4).Value = 'TWh'
$ws.Cells.Item(11,5).Value = 'GW'
$ws.Cells.Item(11,6).Value = 'ANNUAL'

$ws.Cells.Item(12,1).Value = 'ELE'
$ws.Cells.Item(12,2).Value = 'EN_WON_22_c04_CHE'
$ws.Cells.Item(12,3).Value = 'Wind Onshore - CF Class-22 Cost Class-c04 - Switzerland'
$ws.Cells.Item(12,4).Value = 'TWh'
$ws.Cells.Item(12,5).Value = 'GW'
$ws.Cells.Item(12,6).Value = 'ANNUAL'

$ws.Cells.Item(13,1).Value = 'ELE'
$ws.Cells.Item(13,2).Value = 'EN_WON_23_c04_CHE'
$ws.Cells.Item(13,3).Value = 'Wind Onshore - CF Class-23 Cost Class-c04 - Switzerland'
$ws.Cells.Item(13,4).Value = 'TWh'
$ws.Cells.Item(13,5).Value = 'GW'
$ws.Cells.Item(13,6).Value = 'ANNUAL'

$ws.Cells.Item(14,1).Value = 'ELE'
$ws.Cells.Item(14,2).Value = 'EN_WON_24_c02_CHE'
$ws.Cells.Item(14,3).Value = 'Wind Onshore - CF Class-24 Cost Class-c02 - Switzerland'
$ws.Cells.Item(14,4).Value = 'TWh'
$ws.Cells.Item(14,5).Value = 'GW'
$ws.Cells.Item(14,6).Value = 'ANNUAL'

$ws.Cells.Item(15,1).Value = 'ELE'
$ws.Cells.Item(15,2).Value = 'EN_WON_25_c04_CHE'
$ws.Cells.Item(15,3).Value = 'Wind Onshore - CF Class-25 Cost Class-c04 - Switzerland'
$ws.Cells.Item(15,4).Value = 'TWh'
$ws.Cells.Item(15,5).Value = 'GW'
$ws.Cells.Item(15,6).Value = 'ANNUAL'

$ws.Cells.Item(16,1).Value = 'ELE'
$ws.Cells.Item(16,2).Value = 'EN_Hydro_CHE-1'
$ws.Cells.Item(16,3).Value = 'New Hydro Potential - Switzerland - Step 1'
$ws.Cells.Item(16,4).Value = 'PJ'
$ws.Cells.Item(16,5).Value = 'GW'
$ws.Cells.Item(16,6).Value = 'DAYNITE'

$ws.Cells.Item(17,1).Value = 'ELE'
$ws.Cells.Item(17,2).Value = 'EN_Hydro_CHE-2'
$ws.Cells.Item(17,3).Value = 'New Hydro Potential - Switzerland - Step 2'
$ws.Cells.Item(17,4).Value = 'PJ'
$ws.Cells.Item(17,5).Value = 'GW'
$ws.Cells.Item(17,6).Value = 'DAYNITE'

$ws.Cells.Item(18,1).Value = 'ELE'
$ws.Cells.Item(18,2).Value = 'EN_Hydro_CHE-3'
$ws.Cells.Item(18,3).Value = 'New Hydro Potential - Switzerland - Step 3'
$ws.Cells.Item(18,4).Value = 'PJ'
$ws.Cells.Item(18,5).Value = 'GW'
$ws.Cells.Item(18,6).Value = 'DAYNITE'

# Columns I-N (process, Comm-IN, Comm-OUT, CAP_BND, INVCOST~USD21_alt, AF~FX)
$ws.Cells.Item(2,9).Value = 'process'
$ws.Cells.Item(2,10).Value = 'Comm-IN'
$ws.Cells.Item(2,11).Value = 'Comm-OUT'
$ws.Cells.Item(2,12).Value = 'CAP_BND'
$ws.Cells.Item(2,13).Value = 'INVCOST~USD21_alt'
$ws.Cells.Item(2,14).Value = 'AF~FX'

$ws.Cells.Item(3,9).Value = 'EN_SPV_13_c02_CHE'
$ws.Cells.Item(3,10).Value = 'solar'
$ws.Cells.Item(3,11).Value = 'ELC_Sol-CHE'
$ws.Cells.Item(3,12).Value = 22.656000000000006
$ws.Cells.Item(3,13).Value = 88.95499217864773
$ws.Cells.Item(3,14).Value = 0.13052121954449153

$ws.Cells.Item(4,9).Value = 'EN_SPV_13_c03_CHE'
$ws.Cells.Item(4,10).Value = 'solar'
$ws.Cells.Item(4,11).Value = 'ELC_Sol-CHE'
$ws.Cells.Item(4,12).Value = 10.22625
$ws.Cells.Item(4,13).Value = 101.20899620325736
$ws.Cells.Item(4,14).Value = 0.1287234323432343

$ws.Cells.Item(5,9).Value = 'EN_SPV_14_c02_CHE'
$ws.Cells.Item(5,10).Value = 'solar'
$ws.Cells.Item(5,11).Value = 'ELC_Sol-CHE'
$ws.Cells.Item(5,12).Value = 9.573749999999999
$ws.Cells.Item(5,13).Value = 88.95499217864773
$ws.Cells.Item(5,14).Value = 0.13908789659224444

$ws.Cells.Item(6,9).Value = 'EN_SPV_15_c02_CHE'
$ws.Cells.Item(6,10).Value = 'solar'
$ws.Cells.Item(6,11).Value = 'ELC_Sol-CHE'
$ws.Cells.Item(6,12).Value = 0.5009999999999999
$ws.Cells.Item(6,13).Value = 88.95499217864773
$ws.Cells.Item(6,14).Value = 0.14621556886227546

$ws.Cells.Item(7,9).Value = 'EN_WON_18_c04_CHE'
$ws.Cells.Item(7,10).Value = 'wind'
$ws.Cells.Item(7,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(7,12).Value = 6.7065
$ws.Cells.Item(7,13).Value = 145.23264029166972
$ws.Cells.Item(7,14).Value = 0.1782639230597182

$ws.Cells.Item(8,9).Value = 'EN_WON_20_c03_CHE'
$ws.Cells.Item(8,10).Value = 'wind'
$ws.Cells.Item(8,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(8,12).Value = 4.658249999999999
$ws.Cells.Item(8,13).Value = 105.74751621237203
$ws.Cells.Item(8,14).Value = 0.19600000000000004

$ws.Cells.Item(9,9).Value = 'EN_WON_20_c04_CHE'
$ws.Cells.Item(9,10).Value = 'wind'
$ws.Cells.Item(9,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(9,12).Value = 0.3052499999999999
$ws.Cells.Item(9,13).Value = 145.23264029166972
$ws.Cells.Item(9,14).Value = 0.196

$ws.Cells.Item(10,9).Value = 'EN_WON_22_c02_CHE'
$ws.Cells.Item(10,10).Value = 'wind'
$ws.Cells.Item(10,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(10,12).Value = 7.6979999999999995
$ws.Cells.Item(10,13).Value = 80.7856561622413
$ws.Cells.Item(10,14).Value = 0.22300000000000003

$ws.Cells.Item(11,9).Value = 'EN_WON_22_c03_CHE'
$ws.Cells.Item(11,10).Value = 'wind'
$ws.Cells.Item(11,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(11,12).Value = 14.13675
$ws.Cells.Item(11,13).Value = 105.74751621237203
$ws.Cells.Item(11,14).Value = 0.21952581038781904

$ws.Cells.Item(12,9).Value = 'EN_WON_22_c04_CHE'
$ws.Cells.Item(12,10).Value = 'wind'
$ws.Cells.Item(12,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(12,12).Value = 0.46799999999999997
$ws.Cells.Item(12,13).Value = 145.23264029166972
$ws.Cells.Item(12,14).Value = 0.221

$ws.Cells.Item(13,9).Value = 'EN_WON_23_c04_CHE'
$ws.Cells.Item(13,10).Value = 'wind'
$ws.Cells.Item(13,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(13,12).Value = 2.4682500000000003
$ws.Cells.Item(13,13).Value = 145.23264029166972
$ws.Cells.Item(13,14).Value = 0.22986782133090247

$ws.Cells.Item(14,9).Value = 'EN_WON_24_c02_CHE'
$ws.Cells.Item(14,10).Value = 'wind'
$ws.Cells.Item(14,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(14,12).Value = 7.5885
$ws.Cells.Item(14,13).Value = 80.7856561622413
$ws.Cells.Item(14,14).Value = 0.242

$ws.Cells.Item(15,9).Value = 'EN_WON_25_c04_CHE'
$ws.Cells.Item(15,10).Value = 'wind'
$ws.Cells.Item(15,11).Value = 'ELC_Win-CHE'
$ws.Cells.Item(15,12).Value = 0.4484999999999999
$ws.Cells.Item(15,13).Value = 145.23264029166972
$ws.Cells.Item(15,14).Value = 0.246

$ws.Cells.Item(16,9).Value = 'EN_Hydro_CHE-1'
$ws.Cells.Item(16,10).Value = 'hydro'
$ws.Cells.Item(16,11).Value = 'ELC'
$ws.Cells.Item(16,12).Value = 2.1950000000000003

$ws.Cells.Item(17,9).Value = 'EN_Hydro_CHE-2'
$ws.Cells.Item(17,10).Value = 'hydro'
$ws.Cells.Item(17,11).Value = 'ELC'
$ws.Cells.Item(17,12).Value = 6.92

$ws.Cells.Item(18,9).Value = 'EN_Hydro_CHE-3'
$ws.Cells.Item(18,10).Value = 'hydro'
$ws.Cells.Item(18,11).Value = 'ELC'
$ws.Cells.Item(18,12).Value = 13.84
